# Add the new GDP contribution code/label to the INDICATOR sheet and
# switch the active sheet/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INDICATOR")

# New row of data: RLGDPCNT / Real GDP Contribution
$ws.Range("A6").Value = "RLGDPCNT"
$ws.Range("B6").Value = "Real GDP Contribution"

# Column A widened to fit the new, longer code value.
$ws.Columns.Item(1).ColumnWidth = 9.5

# Make INDICATOR the active sheet with the given selection.
$ws.Activate()
$ws.Range("F10").Select()

$wb.Save()
